$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.265.94"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "1.882.21"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'316.45"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.4323"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").Value = "'0.3709"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").Value = "'0.07434"
$ws.Range("D10").Value = "'0.8924"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'21.22"
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("D12").Value = "1.873.47"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "'5.453"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "'6.644"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'0.06986"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "'81.36"
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("D18").Value = "'0.000009129"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'15.67"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "28.360.92"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "'5.098"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "'11.03"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").Value = "2.156.31"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").Value = "'1.995"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "'153.76"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'18.77"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "'5.454"
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("D29").Value = "'118.07"
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("D30").Value = "'1.910"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'0.08989"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").Value = "'0.8023"
$ws.Range("E32").Value = "  +5.08%  "
$ws.Range("D33").Value = "'4.696"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").Value = "'1.185"
$ws.Range("E34").Value = "  +7.59%  "
$ws.Range("D35").Value = "'3.010"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'1.133"
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'0.05498"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "'0.01974"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "'2.899"
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").Value = "'0.1702"
$ws.Range("E41").Value = "  +2.70%  "
$ws.Range("D42").Value = "'0.5193"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").Value = "'6.905"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'8.628"
$ws.Range("E44").Value = "  +4.37%  "
$ws.Range("D45").Value = "'10.62"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("D46").Value = "'0.06618"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'0.4786"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'106.18"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'1.668"
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("D51").Value = "'1.888"
$ws.Range("E51").Value = "  +7.54%  "
